$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its existing text formatting (thousands-dot
# separators, fixed decimal places, etc.) instead of Excel auto-converting
# the new values to numbers when assigned via .Value
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "64.567.93"
$ws.Range("E2").Value = "  +1.23%  "

# Row 3
$ws.Range("D3").Value = "3.088.46"
$ws.Range("E3").Value = "  +0.34%  "

# Row 4
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.24%  "

# Row 5
$ws.Range("D5").Value = "559.11"
$ws.Range("E5").Value = "  +1.10%  "

# Row 6
$ws.Range("D6").Value = "146.21"
$ws.Range("E6").Value = "  +5.46%  "

# Row 7
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.28%  "

# Row 8
$ws.Range("D8").Value = "3.089.20"
$ws.Range("E8").Value = "  +0.59%  "

# Row 9
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -0.38%  "

# Row 10
$ws.Range("D10").Value = "6.39"
$ws.Range("E10").Value = "  +3.13%  "

# Row 11
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  +1.00%  "

# Row 12
$ws.Range("D12").Value = "0.474"
$ws.Range("E12").Value = "  +4.07%  "

# Row 13
$ws.Range("D13").Value = "0.0000230"
$ws.Range("E13").Value = "  +1.21%  "

# Row 14
$ws.Range("D14").Value = "35.49"
$ws.Range("E14").Value = "  +1.58%  "

# Row 15
$ws.Range("D15").Value = "3.587.10"
$ws.Range("E15").Value = "  +0.28%  "

# Row 16
$ws.Range("D16").Value = "64.486.76"
$ws.Range("E16").Value = "  +1.10%  "

# Row 17
$ws.Range("D17").Value = "3.076.60"
$ws.Range("E17").Value = "  -0.10%  "

# Row 18
$ws.Range("E18").Value = "  +0.70%  "

# Row 19
$ws.Range("D19").Value = "6.80"
$ws.Range("E19").Value = "  +0.86%  "

# Row 20
$ws.Range("D20").Value = "477.95"
$ws.Range("E20").Value = "  -1.69%  "

# Row 21
$ws.Range("D21").Value = "14.00"
$ws.Range("E21").Value = "  +3.42%  "

# Row 22
$ws.Range("D22").Value = "0.685"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("D23").Value = "7.55"
$ws.Range("E23").Value = "  +4.96%  "

# Row 24
$ws.Range("D24").Value = "13.59"
$ws.Range("E24").Value = "  +8.30%  "

# Row 25
$ws.Range("D25").Value = "81.13"
$ws.Range("E25").Value = "  -0.48%  "

# Row 26
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  +1.12%  "

# Row 28
$ws.Range("D28").Value = "8.21"
$ws.Range("E28").Value = "  +2.41%  "

# Row 29
$ws.Range("D29").Value = "2.06"
$ws.Range("E29").Value = "  +3.42%  "

# Row 30
$ws.Range("D30").Value = "0.997"
$ws.Range("E30").Value = "  -0.35%  "

# Row 31
$ws.Range("D31").Value = "26.15"
$ws.Range("E31").Value = "  +0.37%  "

# Row 32
$ws.Range("E32").Value = "  +0.29%  "

# Row 33
$ws.Range("D33").Value = "2.50"
$ws.Range("E33").Value = "  +3.07%  "

# Row 34
$ws.Range("D34").Value = "5.62"
$ws.Range("E34").Value = "  -3.29%  "

# Row 35
$ws.Range("D35").Value = "6.17"
$ws.Range("E35").Value = "  +3.06%  "

# Row 36
$ws.Range("D36").Value = "54.87"
$ws.Range("E36").Value = "  -1.54%  "

# Row 37
$ws.Range("D37").Value = "467.85"
$ws.Range("E37").Value = "  -0.85%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").Value = "3.01"
$ws.Range("E38").Value = "  +17.02%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "0.0836"
$ws.Range("E39").Value = "  +2.17%  "

# Row 40
$ws.Range("D40").Value = "0.0408"
$ws.Range("E40").Value = "  +2.86%  "

# Row 41
$ws.Range("D41").Value = "2.970.11"
$ws.Range("E41").Value = "  -6.91%  "

# Row 42
$ws.Range("D42").Value = "8.28"
$ws.Range("E42").Value = "  +0.64%  "

# Row 43
$ws.Range("E43").Value = "  -4.40%  "

# Row 44
$ws.Range("D44").Value = "28.65"
$ws.Range("E44").Value = "  +2.33%  "

# Row 45
$ws.Range("D45").Value = "0.260"
$ws.Range("E45").Value = "  +2.86%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.15"
$ws.Range("E46").Value = "  +5.49%  "

# Row 47
$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48
$ws.Range("E48").Value = "  +2.18%  "

# Row 49
$ws.Range("D49").Value = "0.0₃0524"
$ws.Range("E49").Value = "  +1.65%  "

# Row 50
$ws.Range("D50").Value = "118.23"
$ws.Range("E50").Value = "  +1.39%  "

# Row 51
$ws.Range("D51").Value = "2.08"
$ws.Range("E51").Value = "  +0.36%  "
